# add: content song hindia - topik
# Replaces the 3 placeholder song rows (Ascence / Cartoon / NCS samples) with the
# full "Hindia - Topik" album track list (12 songs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new rows (5-13) need the same "no-hyperlink" style that column A already
# carries on rows 2-4 (s="1"). Copy that formatting down before writing values.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5:A13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# songImage, songMusic, songName, artistName, albumName
$ws.Range("A2").Value2 = "ic_hindia_cover"
$ws.Range("B2").Value2 = "apapun_yang_terjadi"
$ws.Range("C2").Value2 = "Apapun Yang Terjadi"
$ws.Range("D2").Value2 = "Hindia"
$ws.Range("E2").Value2 = "Topik"

$ws.Range("A3").Value2 = "ic_hindia_cover"
$ws.Range("B3").Value2 = "belum_tidur_feat_sal_priadi"
$ws.Range("C3").Value2 = "Belum Tidur (feat. Sal Priadi)"
$ws.Range("D3").Value2 = "Hindia"
$ws.Range("E3").Value2 = "Topik"

$ws.Range("A4").Value2 = "ic_hindia_cover"
$ws.Range("B4").Value2 = "besok_mungkin_kita_sampai"
$ws.Range("C4").Value2 = "Besok Mungkin Kita Sampai"
$ws.Range("D4").Value2 = "Hindia"
$ws.Range("E4").Value2 = "Topik"

$ws.Range("A5").Value2 = "ic_hindia_cover"
$ws.Range("B5").Value2 = "dehidrasi"
$ws.Range("C5").Value2 = "Dehidrasi (feat. Petra Sihombing)"
$ws.Range("D5").Value2 = "Hindia"
$ws.Range("E5").Value2 = "Topik"

$ws.Range("A6").Value2 = "ic_hindia_cover"
$ws.Range("B6").Value2 = "evakuasi"
$ws.Range("C6").Value2 = "Evakuasi"
$ws.Range("D6").Value2 = "Hindia"
$ws.Range("E6").Value2 = "Topik"

$ws.Range("A7").Value2 = "ic_hindia_cover"
$ws.Range("B7").Value2 = "evaluasi"
$ws.Range("C7").Value2 = "Evaluasi"
$ws.Range("D7").Value2 = "Hindia"
$ws.Range("E7").Value2 = "Topik"

$ws.Range("A8").Value2 = "ic_hindia_cover"
$ws.Range("B8").Value2 = "jam_makan_siang"
$ws.Range("C8").Value2 = "Jam Makan Siang (feat. Matter Mos)"
$ws.Range("D8").Value2 = "Hindia"
$ws.Range("E8").Value2 = "Topik"

$ws.Range("A9").Value2 = "ic_hindia_cover"
$ws.Range("B9").Value2 = "mata_air"
$ws.Range("C9").Value2 = "Mata Air (feat. Natasha Udu, Kamga)"
$ws.Range("D9").Value2 = "Hindia"
$ws.Range("E9").Value2 = "Topik"

$ws.Range("A10").Value2 = "ic_hindia_cover"
$ws.Range("B10").Value2 = "membasuh"
$ws.Range("C10").Value2 = "Membasuh (feat. Rara Sekar)"
$ws.Range("D10").Value2 = "Hindia"
$ws.Range("E10").Value2 = "Topik"

$ws.Range("A11").Value2 = "ic_hindia_cover"
$ws.Range("B11").Value2 = "rumah_ke_rumah"
$ws.Range("C11").Value2 = "Rumah Ke Rumah"
$ws.Range("D11").Value2 = "Hindia"
$ws.Range("E11").Value2 = "Topik"

$ws.Range("A12").Value2 = "ic_hindia_cover"
$ws.Range("B12").Value2 = "secukupnya"
$ws.Range("C12").Value2 = "Secukupnya"
$ws.Range("D12").Value2 = "Hindia"
$ws.Range("E12").Value2 = "Topik"

$ws.Range("A13").Value2 = "ic_hindia_cover"
$ws.Range("B13").Value2 = "untuk_apa"
$ws.Range("C13").Value2 = "Untuk Apa - Untuk Apa"
$ws.Range("D13").Value2 = "Hindia"
$ws.Range("E13").Value2 = "Topik"

# Match the final cursor/selection position left by the author's Excel session
$ws.Range("A18").Select() | Out-Null
